# Append the 2024-09-30 requirement rows (87-134) scraped after the last
# existing row (86, dated 2024-09-29) to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# Columns B (VendorPhone) and C (Date) must stay TEXT: some phone numbers
# are purely numeric (e.g. "17866280721", "125.36132370") and every Date
# value is the literal string "2024-09-30" -- without forcing a text
# number format up front, Excel would silently reinterpret these as a
# number / a real date serial when .Value is assigned.
$textRange = $ws.Range("B87:C134")
$textRange.NumberFormat = "@"

$ws.Range("A87").Value = '"Isaac, Centraprise" <isaac@centraprise.com>'
$ws.Range("B87").Value = '469-923-8111'
$ws.Range("C87").Value = '2024-09-30'
$ws.Range("D87").Value = 'Java Developer at NYC NY (Day One Onsite) Need Ex Cognizant OR Ex Amex'

$ws.Range("A88").Value = '"Mehak Pandey, Simplify Software Experts" <mehak@simplifysoftwareexperts.com>'
$ws.Range("B88").Value = '(201)2855932'
$ws.Range("C88").Value = '2024-09-30'
$ws.Range("D88").Value = 'CORE JAVA BACKEND DEVELOPER'

$ws.Range("A89").Value = '"Raman Arora, Zealhire" <raman@zealhire.com>'
$ws.Range("B89").Value = '(332) 2878468'
$ws.Range("C89").Value = '2024-09-30'
$ws.Range("D89").Value = 'IT Oracle Project Managers :: hybrid role :: Local Only'

$ws.Range("A90").Value = 'Vishwas Vishwas <vishwasvishwas622@gmail.com>'
$ws.Range("B90").Value = '55867-92544'
$ws.Range("C90").Value = '2024-09-30'
$ws.Range("D90").Value = 'Sr Golang Developers || Rosemont, IL/Onsite)'

$ws.Range("A91").Value = 'amit <amitcorp2corp@gmail.com>'
$ws.Range("B91").Value = '848-666-1490'
$ws.Range("C91").Value = '2024-09-30'
$ws.Range("D91").Value = 'Opening for Oracle Supplier Data Hub (Oracle SDH)'

$ws.Range("A92").Value = 'amit <amitcorp2corp@gmail.com>'
$ws.Range("B92").Value = '848-666-1490'
$ws.Range("C92").Value = '2024-09-30'
$ws.Range("D92").Value = '{C2C-W2 -Jobs} Opening for Oracle Supplier Data Hub (Oracle SDH)'

$ws.Range("A93").Value = 'Pankaj Chauhan <pchauhan@accroid.com>'
$ws.Range("B93").Value = '(302)-4851559'
$ws.Range("C93").Value = '2024-09-30'
$ws.Range("D93").Value = 'Oracle EBS IT Project Manager || Minneapolis, MN (On-Site) relocation' + $nl + ' works  || Visa: USC/GC'

$ws.Range("A94").Value = 'ajay.sai@dprsolutionsinc.com'
$ws.Range("B94").Value = '5714632218, (603) 4860938'
$ws.Range("C94").Value = '2024-09-30'
$ws.Range("D94").Value = 'Resume Submission Java Full Stack Developer'

$ws.Range("A95").Value = 'Vineeth Damarla <ron@americanitsystems.com>'
$ws.Range("B95").Value = '125.36132370'
$ws.Range("C95").Value = '2024-09-30'
$ws.Range("D95").Value = 'Java Technical Developer | Berkeley Heights, NJ 100% onsite from' + $nl + ' day-1)'

$ws.Range("A96").Value = 'Vineeth Damarla <ron@americanitsystems.com>'
$ws.Range("B96").Value = '125.36132370'
$ws.Range("C96").Value = '2024-09-30'
$ws.Range("D96").Value = 'Java Technical Developer | Berkeley Heights, NJ (100% onsite from' + $nl + ' day-1)'

$ws.Range("A97").Value = '"Shivam Tayal, Quantum World Technologies Inc" <shivam.tayal@quantumworldit.com>'
$ws.Range("B97").Value = '+1 805-6784659'
$ws.Range("C97").Value = '2024-09-30'
$ws.Range("D97").Value = 'Hiring Oracle Sales Cloud SDH Consultant, Chicago IL'

$ws.Range("A98").Value = '"sachin, Atlas Cloud, solution" <sachin@acloudinc.com>'
$ws.Range("B98").Value = '+1 4195959070'
$ws.Range("C98").Value = '2024-09-30'
$ws.Range("D98").Value = 'Microsoft Dynamics D365'

$ws.Range("A99").Value = '"Vignesh Ramakrishnan, Teamware solution" <vignesh.r@twsol.com>'
$ws.Range("B99").Value = '+12145449254'
$ws.Range("C99").Value = '2024-09-30'
$ws.Range("D99").Value = 'Need-Oracle Supplier Data Hub consultant Chicago, IL'

$ws.Range("A100").Value = 'Jaffer Shahul H <shahul.h@diamondpick.com>'
$ws.Range("B100").Value = '34843418293, 6042538933'
$ws.Range("C100").Value = '2024-09-30'
$ws.Range("D100").Value = 'Immediate Hire - Microsoft SQL Server (DBA - Architect/ Consultant) -' + $nl + ' Reading, PA(Hybrid) - Only Local'

$ws.Range("A101").Value = '"irfan shiak, Agile Enterprise Solutions Inc." <irfan_shaik@aesincus.com>'
$ws.Range("B101").Value = '647-375-2228, 972-440-0069'
$ws.Range("C101").Value = '2024-09-30'
$ws.Range("D101").Value = 'SingleStore DBA with Python Programming |Austin, TX| Onsite'

$ws.Range("A102").Value = '"Suzanne Rogers, Concord IT Systems" <suzanne.rogers@concorditsystems.com>'
$ws.Range("B102").Value = '17866280721'
$ws.Range("C102").Value = '2024-09-30'
$ws.Range("D102").Value = 'Immediate hiring for Oracle NetSuite Functional Consultant - Remote'

$ws.Range("A103").Value = '"Satish Kumar, Donato Technologies INC" <satish@donatotech.net>'
$ws.Range("B103").Value = '(469) 9299409, 469-533-0333'
$ws.Range("C103").Value = '2024-09-30'
$ws.Range("D103").Value = 'Looking for JPC - 2122 - Oracle IAM Architect -Latham, New York (Onsite)- Contract Job- Sligo'

$ws.Range("A104").Value = 'Surya Hemanth <hemanth@brillius.com>'
$ws.Range("B104").Value = '+1 510-4008680'
$ws.Range("C104").Value = '2024-09-30'
$ws.Range("D104").Value = 'JAVA DEVELOPER WITH ML Exp :: Bellevue, WA'

$ws.Range("A105").Value = '"Suraj Prashar, Pivotal Technologies" <suraj.prashar@pivotal-technologies.com>'
$ws.Range("B105").Value = '(703) 5708775'
$ws.Range("C105").Value = '2024-09-30'
$ws.Range("D105").Value = 'Urgent Requirement || SQL DBA || Remote|| 6+Months'

$ws.Range("A106").Value = '"Pallavi, Yochana" <pallavi@yochana.com>'
$ws.Range("B106").Value = '2482373189'
$ws.Range("C106").Value = '2024-09-30'
$ws.Range("D106").Value = 'Please share local candidates-Sr Golang Developer  Cupertino, CA (Onsite)-Job Description'

$ws.Range("A107").Value = 'Venkat G <venkat.g@stiorg.com>'
$ws.Range("B107").Value = '+1 (740) 3273431, +1 (609) 9983429'
$ws.Range("C107").Value = '2024-09-30'
$ws.Range("D107").Value = 'RE: Sr. Java Lead Consultant available to take new project and open' + $nl + ' to relocate.'

$ws.Range("A108").Value = '"Ishika, Thothit LLC" <ishika@thothit.net>'
$ws.Range("B108").Value = '352-614-4461'
$ws.Range("C108").Value = '2024-09-30'
$ws.Range("D108").Value = 'General Database Admin(this is for a DBA not developer)||  Alpharetta, GA(need local with DL)'

$ws.Range("A109").Value = '"Riya Kori, Gtech LLC" <rkori@greattechglobal.com>'
$ws.Range("B109").Value = '469-527-2150'
$ws.Range("C109").Value = '2024-09-30'
$ws.Range("D109").Value = 'AWS JAVA Engineer || Onsite in Seattle WA || No H1b || Rate is $50 || Ned Aws Active Certification'

$ws.Range("A110").Value = '"Shanu, Parmesoft" <shanu.f@parmesoft.com>'
$ws.Range("B110").Value = '972-402-5580, 289-652-1056'
$ws.Range("C110").Value = '2024-09-30'
$ws.Range("D110").Value = 'Immediate Hire -Oracle EBS PLSQL Developer-Philadelphia, PA(Hybrid)- 12+ plus years of experience'

$ws.Range("A111").Value = '"vaibhav kumar, VBeyond Corp" <vaibhavk@vbeyond.com>'
$ws.Range("B111").Value = '+1-9086334110'
$ws.Range("C111").Value = '2024-09-30'
$ws.Range("D111").Value = 'Sr. Oracle PL-SQL Developer in Chicago, IL (Onsite)'

$ws.Range("A112").Value = 'Murthy Medisetti <murthy@pristineresource.com>'
$ws.Range("B112").Value = '860-515-8853, 204565678100'
$ws.Range("C112").Value = '2024-09-30'
$ws.Range("D112").Value = 'Techno-Functional Architect (Oracle Revenue Management Cloud' + $nl + ' Services)'

$ws.Range("A113").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B113").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C113").Value = '2024-09-30'
$ws.Range("D113").Value = 'Re: Nee Power BI with java and security for remote'

$ws.Range("A114").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B114").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C114").Value = '2024-09-30'
$ws.Range("D114").Value = 'Re: Nee Power BI with java and security for remote'

$ws.Range("A115").Value = '"ayush, istaff" <ayush@istaffx.com>'
$ws.Range("B115").Value = '575-236-4255'
$ws.Range("C115").Value = '2024-09-30'
$ws.Range("D115").Value = 'Java Developer Hybrid role in TX with Selenium Experience Need locals Only in Austin , TX   No h1B'

$ws.Range("A116").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B116").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C116").Value = '2024-09-30'
$ws.Range("D116").Value = 'Re: Nee Power BI with Java and security for remote'

$ws.Range("A117").Value = '"Suman Bala, HMG America" <sbala@hmgamerica.com>'
$ws.Range("B117").Value = '7327905647, 7327905001'
$ws.Range("C117").Value = '2024-09-30'
$ws.Range("D117").Value = '(Onsite)URGENT REQUIREMENT- SQL Server DBA in Plano TX- onsite'

$ws.Range("A118").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B118").Value = '(505) 2125488, 467914261938, (916) 8581390'
$ws.Range("C118").Value = '2024-09-30'
$ws.Range("D118").Value = 'Re: Need Power BI with Java and security for remote'

$ws.Range("A119").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B119").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C119").Value = '2024-09-30'
$ws.Range("D119").Value = 'Re: Need Power BI with Java(Must) and security for remote'

$ws.Range("A120").Value = '<prince.sharma@applabsystems.com>'
$ws.Range("B120").Value = '2463466860, 2421177820, 609-766-0112'
$ws.Range("C120").Value = '2024-09-30'
$ws.Range("D120").Value = 'Oracle Integration Cloud (OIC) || Atlanta, GA; Frisco, TX; Bellevue, WA (Onsite)'

$ws.Range("A121").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B121").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C121").Value = '2024-09-30'
$ws.Range("D121").Value = 'Re: Nee Power BI with java and security for remote'

$ws.Range("A122").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B122").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C122").Value = '2024-09-30'
$ws.Range("D122").Value = 'Re: Nee Power BI with java and security for remote'

$ws.Range("A123").Value = 'Ajay Immadisetty <ajay.immadisetty@sitacorp.com>'
$ws.Range("B123").Value = '+1 732 9067806'
$ws.Range("C123").Value = '2024-09-30'
$ws.Range("D123").Value = 'Oracle HCM Technical Lead/Manager'

$ws.Range("A124").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B124").Value = '(505) 2125488, 42-9833-474986, (916) 8581390'
$ws.Range("C124").Value = '2024-09-30'
$ws.Range("D124").Value = 'Re: Nee Power BI with java and security for remote'

$ws.Range("A125").Value = '"Vikas Chaudhary, DMS VISIONS INC" <vikas@dmsvisions.com>'
$ws.Range("B125").Value = '972-645-0989'
$ws.Range("C125").Value = '2024-09-30'
$ws.Range("D125").Value = 'Urgent Hiring || Java Developer || Columbus OH || Hybrid || USC or GC Only || 06 Months C2H || 2nd Round IN person || Local candidates Only Interview'

$ws.Range("A126").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B126").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C126").Value = '2024-09-30'
$ws.Range("D126").Value = 'Re: Nee Power BI with java and security for remote'

$ws.Range("A127").Value = 'Vijay <vijay@jnjsolutions.com>'
$ws.Range("B127").Value = '(505) 2125488, (916) 8581390'
$ws.Range("C127").Value = '2024-09-30'
$ws.Range("D127").Value = 'Re: Nee Power BI with java and security for remote'

$ws.Range("A128").Value = '"Srikanth, iTech US" <srikanth.v@itechus.net>'
$ws.Range("B128").Value = '802 798 8941'
$ws.Range("C128").Value = '2024-09-30'
$ws.Range("D128").Value = 'Oracle PL | SQL Developer with Informatica CDI &amp; IICS Exp - Maricopa, AZ(Onsite) - MST | PST Candidates'

$ws.Range("A129").Value = '"Gangadar Reddy, Centraprise" <edula.gangadar@centraprise.com>'
$ws.Range("B129").Value = '469-639-0369'
$ws.Range("C129").Value = '2024-09-30'
$ws.Range("D129").Value = 'Hiring for Java  Engineer II - New York, NY - Day 1 Onsite - AMEX-cognizant'

$ws.Range("A130").Value = '"Narsimha Kathi, iTech US" <narsimha.k@itechus.net>'
$ws.Range("B130").Value = '802-383-1500, 802-735-0270'
$ws.Range("C130").Value = '2024-09-30'
$ws.Range("D130").Value = 'Urgent Requirement for Oracle PLSQL Developer - Phoenix, AZ-Maricopa, AZ - Hybrid Locals or Nearby'

$ws.Range("A131").Value = 'Vishwas Vishwas <vishwasvishwas622@gmail.com>'
$ws.Range("B131").Value = '+1 9105576339'
$ws.Range("C131").Value = '2024-09-30'
$ws.Range("D131").Value = 'Re: Sr Golang Developers || Rosemont, IL/Onsite)'

$ws.Range("A132").Value = 'Faisal Siddiqui <thereisnodifferncebet123@gmail.com>'
$ws.Range("B132").Value = '+1 9105576339'
$ws.Range("C132").Value = '2024-09-30'
$ws.Range("D132").Value = 'Re: Fwd: Full Stack Java AWS Developer || Seattle, WA (To go onsite 5' + $nl + ' days a week)'

$ws.Range("A133").Value = '"Rajeev, Tek InspirationsLLC" <rajeev.kharwar@tekinspirations.com>'
$ws.Range("B133").Value = '7525894499, 469-393-0216'
$ws.Range("C133").Value = '2024-09-30'
$ws.Range("D133").Value = 'AWS PYTHON GOLANG DEVELOPER'

$ws.Range("A134").Value = '"Yashwant Singh, DMS Visions Inc" <yashwant@dmsvisions.com>'
$ws.Range("B134").Value = '972-645-5050'
$ws.Range("C134").Value = '2024-09-30'
$ws.Range("D134").Value = 'Java Developer with CMS and Retail  (Hybrid) Columbus, OH || Locals only || USC &amp; GC only'

# The temporary "@" text format left a style index on B87:C134; reset it
# back to the default/Normal style so the new cells match the rest of the
# sheet (the pre-existing data rows carry no explicit cell style either).
$textRange.Style = "Normal"
